$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 49130.906
$ws.Range("J17").Value = 49130.906
$ws.Range("L17").Value = 147392.718
$ws.Range("N17").Value = -147728.718
$ws.Range("H32").Value = 5554.273
$ws.Range("J32").Value = 6899.625
$ws.Range("L32").Value = 6899.625
$ws.Range("N32").Value = -7551.625
$ws.Range("H106").Value = 9756.546
$ws.Range("I106").Value = 9756.546
$ws.Range("K106").Value = 9756.546
$ws.Range("M106").Value = -9125.546
$ws.Range("H130").Value = 85136.5
$ws.Range("J130").Value = 85136.5
$ws.Range("L130").Value = 85136.5
$ws.Range("N130").Value = -95176.5
$ws.Range("H132").Value = 1972.4
$ws.Range("I132").Value = 1987.7858
$ws.Range("J132").Value = 1757
$ws.Range("K132").Value = 5963.357400000001
$ws.Range("L132").Value = 5271
$ws.Range("M132").Value = -3433.357400000001
$ws.Range("N132").Value = -10331
$ws.Range("H137").Value = 4382.923
$ws.Range("I137").Value = 4382.923
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 13148.769
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -10598.769
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 11274.56
$ws.Range("I138").Value = 7250.7
$ws.Range("J138").Value = 11721.655
$ws.Range("K138").Value = 21752.1
$ws.Range("L138").Value = 35164.965
$ws.Range("M138").Value = -16612.1
$ws.Range("N138").Value = -45444.965

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1249.75
$ws.Range("I5").Value = 999.6667
$ws.Range("K5").Value = 999.6667
$ws.Range("M5").Value = -887.6667
$ws.Range("H32").Value = 19450.738
$ws.Range("I32").Value = 19450.738
$ws.Range("K32").Value = 19450.738
$ws.Range("M32").Value = -19163.738
$ws.Range("H61").Value = 9549747
$ws.Range("I61").Value = 13337007
$ws.Range("J61").Value = 1434189.4
$ws.Range("K61").Value = 13337007
$ws.Range("L61").Value = 1434189.4
$ws.Range("M61").Value = -13336795
$ws.Range("N61").Value = -1434613.4
$ws.Range("H132").Value = 8335853.5
$ws.Range("I132").Value = 2881.1
$ws.Range("J132").Value = 50000716
$ws.Range("K132").Value = 8643.299999999999
$ws.Range("L132").Value = 150002148
$ws.Range("M132").Value = -6113.299999999999
$ws.Range("N132").Value = -150007208
$ws.Range("H136").Value = 9549747
$ws.Range("I136").Value = 13337007
$ws.Range("J136").Value = 1434189.4
$ws.Range("K136").Value = 40011021
$ws.Range("L136").Value = 4302568.199999999
$ws.Range("M136").Value = -40008471
$ws.Range("N136").Value = -4307668.199999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1249.75
$ws.Range("I4").Value = 999.6667
$ws.Range("K4").Value = 999.6667
$ws.Range("M4").Value = -884.6667
$ws.Range("H75").Value = 10029.571
$ws.Range("I75").Value = 10029.571
$ws.Range("K75").Value = 10029.571
$ws.Range("M75").Value = -9093.571
$ws.Range("H78").Value = 10029.571
$ws.Range("I78").Value = 10029.571
$ws.Range("K78").Value = 30088.713
$ws.Range("M78").Value = -25408.713
$ws.Range("H94").Value = 4056.2354
$ws.Range("I94").Value = 4262.875
$ws.Range("J94").Value = 750
$ws.Range("K94").Value = 4262.875
$ws.Range("L94").Value = 750
$ws.Range("M94").Value = -3811.875
$ws.Range("N94").Value = -1652
$ws.Range("H105").Value = 341195.44
$ws.Range("I105").Value = 430037.3
$ws.Range("K105").Value = 430037.3
$ws.Range("M105").Value = -428290.3
$ws.Range("H134").Value = 6669774
$ws.Range("I134").Value = 2925.7273
$ws.Range("K134").Value = 8777.1819
$ws.Range("M134").Value = -6242.1819

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35719040
$ws.Range("I31").Value = 62503292
$ws.Range("K31").Value = 62503292
$ws.Range("M31").Value = -62502997
$ws.Range("H34").Value = 35719040
$ws.Range("I34").Value = 62503292
$ws.Range("K34").Value = 62503292
$ws.Range("M34").Value = -62503090
$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H59").Value = 87129
$ws.Range("J59").Value = 88999.39999999999
$ws.Range("L59").Value = 88999.39999999999
$ws.Range("N59").Value = -91289.39999999999
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H70").Value = 70000
$ws.Range("J70").Value = 70000
$ws.Range("L70").Value = 70000
$ws.Range("N70").Value = -70630
$ws.Range("H73").Value = 70000
$ws.Range("J73").Value = 70000
$ws.Range("L73").Value = 70000
$ws.Range("N73").Value = -72184
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H122").Value = 6050
$ws.Range("I122").Value = 6050
$ws.Range("K122").Value = 18150
$ws.Range("M122").Value = -15700
$ws.Range("H132").Value = 2627.8635
$ws.Range("I132").Value = 2710.65
$ws.Range("K132").Value = 8131.950000000001
$ws.Range("M132").Value = -5601.950000000001
$ws.Range("H134").Value = 1534.2307
$ws.Range("I134").Value = 1267.7273
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 3803.1819
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -1268.1819
$ws.Range("N134").Value = -14070
$ws.Range("H141").Value = 564195.4
$ws.Range("I141").Value = 20000
$ws.Range("J141").Value = 606056.5600000001
$ws.Range("K141").Value = 20000
$ws.Range("L141").Value = 606056.5600000001
$ws.Range("N141").Value = -616416.5600000001
$ws.Range("M141").Value = -14820

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1282.6666
$ws.Range("J23").Value = 1481.9
$ws.Range("L23").Value = 4445.700000000001
$ws.Range("N23").Value = -4915.700000000001
$ws.Range("H38").Value = 34.576923
$ws.Range("I38").Value = 41.9
$ws.Range("J38").Value = 10.166667
$ws.Range("K38").Value = 125.7
$ws.Range("L38").Value = 30.500001
$ws.Range("M38").Value = 221.3
$ws.Range("N38").Value = -724.500001
$ws.Range("H42").Value = 200002220
$ws.Range("J42").Value = 4000
$ws.Range("L42").Value = 12000
$ws.Range("N42").Value = -13068
$ws.Range("H118").Value = 998.6667
$ws.Range("I118").Value = 998.6667
$ws.Range("K118").Value = 2996.0001
$ws.Range("M118").Value = -1753.0001
$ws.Range("H131").Value = 4337.1
$ws.Range("J131").Value = 3941.8823
$ws.Range("L131").Value = 11825.6469
$ws.Range("N131").Value = -21905.6469
$ws.Range("H134").Value = 10548.5
$ws.Range("I134").Value = 4498.2856
$ws.Range("J134").Value = 24665.666
$ws.Range("K134").Value = 13494.8568
$ws.Range("L134").Value = 73996.99800000001
$ws.Range("M134").Value = -8424.856800000001
$ws.Range("N134").Value = -84136.99800000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 8041.952
$ws.Range("I122").Value = 5944.05
$ws.Range("K122").Value = 17832.15
$ws.Range("M122").Value = -15382.15
$ws.Range("H132").Value = 6367950.5
$ws.Range("I132").Value = 4521.6875
$ws.Range("K132").Value = 13565.0625
$ws.Range("M132").Value = -11035.0625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9114.6
$ws.Range("I40").Value = 5772.522
$ws.Range("J40").Value = 13636.235
$ws.Range("K40").Value = 5772.522
$ws.Range("L40").Value = 13636.235
$ws.Range("M40").Value = -5636.522
$ws.Range("N40").Value = -13908.235
$ws.Range("H46").Value = 1217
$ws.Range("J46").Value = 1204
$ws.Range("L46").Value = 1204
$ws.Range("N46").Value = -1580
$ws.Range("H93").Value = 11116731
$ws.Range("I93").Value = 2000
$ws.Range("J93").Value = 13895414
$ws.Range("K93").Value = 2000
$ws.Range("L93").Value = 13895414
$ws.Range("M93").Value = -752
$ws.Range("N93").Value = -13897910
$ws.Range("H132").Value = 4040.8572
$ws.Range("I132").Value = 2591.7856
$ws.Range("J132").Value = 6939
$ws.Range("K132").Value = 7775.3568
$ws.Range("L132").Value = 20817
$ws.Range("M132").Value = -5245.3568
$ws.Range("N132").Value = -25877
$ws.Range("H133").Value = 89490.25
$ws.Range("J133").Value = 89490.25
$ws.Range("L133").Value = 89490.25
$ws.Range("N133").Value = -94550.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2478.9333
$ws.Range("I100").Value = 1767.625
$ws.Range("K100").Value = 3535.25
$ws.Range("M100").Value = -2994.25
$ws.Range("H122").Value = 2122.4783
$ws.Range("I122").Value = 1968.9546
$ws.Range("K122").Value = 5906.8638
$ws.Range("M122").Value = -3456.8638
